$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.584742307662964
$ws.Range("B1").Value = 1.934595704078674
$ws.Range("C1").Value = 2.090032339096069
$ws.Range("D1").Value = 2.412392616271973
$ws.Range("E1").Value = 3.217165470123291
